$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC, @@ -4776,25 +4776,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 4135
$ws.Range("J82").Value = 8000
$ws.Range("L82").Value = 24000
$ws.Range("N82").Value = -24812

# Hunk 1: sheet ALC, @@ -4929,25 +4929,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 4135
$ws.Range("J85").Value = 8000
$ws.Range("L85").Value = 24000
$ws.Range("N85").Value = -26808

# Hunk 2: sheet ALC, @@ -7142,22 +7142,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 4387089
$ws.Range("I129").Value = 35715440
$ws.Range("K129").Value = 107146320
$ws.Range("M129").Value = -107141320

# Hunk 3: sheet ALC, @@ -7292,22 +7292,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4350379.5
$ws.Range("I132").Value = 5002238.5
$ws.Range("K132").Value = 15006715.5
$ws.Range("M132").Value = -15004185.5

# Hunk 4: sheet ALC, @@ -7442,25 +7442,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 811.6458
$ws.Range("I135").Value = 561.275
$ws.Range("J135").Value = 2063.5
$ws.Range("K135").Value = 5051.474999999999
$ws.Range("L135").Value = 18571.5
$ws.Range("M135").Value = -2516.474999999999
$ws.Range("N135").Value = -23641.5

# Hunk 5: sheet ALC, @@ -7543,25 +7543,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2484.2156
$ws.Range("I137").Value = 2370.175
$ws.Range("J137").Value = 2898.9092
$ws.Range("K137").Value = 7110.525000000001
$ws.Range("L137").Value = 8696.7276
$ws.Range("M137").Value = -4560.525000000001
$ws.Range("N137").Value = -13796.7276

# Hunk 6: sheet ARM, @@ -10025,25 +10025,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1395.2903
$ws.Range("I45").Value = 1068
$ws.Range("J45").Value = 4450
$ws.Range("K45").Value = 1068
$ws.Range("L45").Value = 4450
$ws.Range("M45").Value = -691
$ws.Range("N45").Value = -5204

# Hunk 7: sheet ARM, @@ -10901,22 +10901,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5866.5
$ws.Range("I63").Value = 3200
$ws.Range("J63").Value = 6399.8
$ws.Range("K63").Value = 3200
$ws.Range("L63").Value = 6399.8
$ws.Range("M63").Value = -2514
$ws.Range("N63").Value = -7771.8

# Hunk 8: sheet ARM, @@ -11045,22 +11048,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 5866.5
$ws.Range("I66").Value = 3200
$ws.Range("J66").Value = 6399.8
$ws.Range("K66").Value = 16000
$ws.Range("L66").Value = 31999
$ws.Range("M66").Value = -12568
$ws.Range("N66").Value = -38863

# Hunk 9: sheet ARM, @@ -11725,22 +11731,22 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 27009
$ws.Range("J80").Value = 27009
$ws.Range("L80").Value = 27009
$ws.Range("N80").Value = -29005

# Hunk 10: sheet ARM, @@ -11872,22 +11878,22 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 27009
$ws.Range("J83").Value = 27009
$ws.Range("L83").Value = 81027
$ws.Range("N83").Value = -91011

# Hunk 11: sheet ARM, @@ -12549,25 +12555,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 822.55
$ws.Range("I97").Value = 816.93335
$ws.Range("J97").Value = 839.4
$ws.Range("K97").Value = 816.93335
$ws.Range("L97").Value = 839.4
$ws.Range("M97").Value = -320.93335
$ws.Range("N97").Value = -1831.4

# Hunk 12: sheet ARM, @@ -12797,22 +12803,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4277.778
$ws.Range("I102").Value = 3562.5
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 3562.5
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -1940.5
$ws.Range("N102").Value = -13244

# Hunk 13: sheet ARM, @@ -14255,25 +14264,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2011.8036
$ws.Range("I132").Value = 1418.875
$ws.Range("J132").Value = 3494.125
$ws.Range("K132").Value = 4256.625
$ws.Range("L132").Value = 10482.375
$ws.Range("M132").Value = -1726.625
$ws.Range("N132").Value = -15542.375

# Hunk 14: sheet BSM, @@ -19916,25 +19925,25 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1478.1628
$ws.Range("I105").Value = 1301.5
$ws.Range("J105").Value = 1748.3529
$ws.Range("K105").Value = 1301.5
$ws.Range("L105").Value = 1748.3529
$ws.Range("M105").Value = 445.5
$ws.Range("N105").Value = -5242.3529

# Hunk 15: sheet BSM, @@ -21227,22 +21236,19 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# Hunk 16: sheet BSM, @@ -21325,25 +21331,25 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1923.7872
$ws.Range("I134").Value = 1311.225
$ws.Range("J134").Value = 5424.143
$ws.Range("K134").Value = 3933.675
$ws.Range("L134").Value = 16272.429
$ws.Range("M134").Value = -1398.675
$ws.Range("N134").Value = -21342.429

# Hunk 17: sheet CRP, @@ -23274,25 +23280,25 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2028.4938
$ws.Range("I31").Value = 1300.7693
$ws.Range("J31").Value = 3333.3794
$ws.Range("K31").Value = 1300.7693
$ws.Range("L31").Value = 3333.3794
$ws.Range("M31").Value = -1005.7693
$ws.Range("N31").Value = -3923.3794

# Hunk 18: sheet CRP, @@ -23424,25 +23430,25 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2028.4938
$ws.Range("I34").Value = 1300.7693
$ws.Range("J34").Value = 3333.3794
$ws.Range("K34").Value = 1300.7693
$ws.Range("L34").Value = 3333.3794
$ws.Range("M34").Value = -1098.7693
$ws.Range("N34").Value = -3737.3794

# Hunk 19: sheet CRP, @@ -28247,25 +28253,25 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2666.3462
$ws.Range("I132").Value = 1949.8823
$ws.Range("J132").Value = 4019.6667
$ws.Range("K132").Value = 5849.6469
$ws.Range("L132").Value = 12059.0001
$ws.Range("M132").Value = -3319.6469
$ws.Range("N132").Value = -17119.0001

# Hunk 20: sheet CUL, @@ -35777,25 +35783,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3865.6667
$ws.Range("J138").Value = 6810
$ws.Range("L138").Value = 20430
$ws.Range("N138").Value = -30710

# Hunk 21: sheet CUL, @@ -35881,22 +35887,22 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 9262525
$ws.Range("I140").Value = 18519308
$ws.Range("K140").Value = 55557924
$ws.Range("M140").Value = -55552744

# Hunk 22: sheet GSM, @@ -38280,22 +38286,22 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 9795.25
$ws.Range("J46").Value = 9795.25
$ws.Range("L46").Value = 9795.25
$ws.Range("N46").Value = -10107.25

# Hunk 23: sheet GSM, @@ -39934,25 +39940,25 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2918.75
$ws.Range("I80").Value = 2835.7144
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 2835.7144
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -1837.7144
$ws.Range("N80").Value = -5496

# Hunk 24: sheet GSM, @@ -40087,25 +40093,25 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2918.75
$ws.Range("I83").Value = 2835.7144
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 14178.572
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -9186.572
$ws.Range("N83").Value = -27484

# Hunk 25: sheet GSM, @@ -42464,22 +42470,22 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2784.2778
$ws.Range("I132").Value = 2318.1765
$ws.Range("K132").Value = 6954.529500000001
$ws.Range("M132").Value = -4424.529500000001

# Hunk 26: sheet LTW, @@ -47522,25 +47528,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2469.9473
$ws.Range("I93").Value = 2176.4
$ws.Range("J93").Value = 3570.75
$ws.Range("K93").Value = 2176.4
$ws.Range("L93").Value = 3570.75
$ws.Range("M93").Value = -928.4000000000001
$ws.Range("N93").Value = -6066.75

# Hunk 27: sheet LTW, @@ -49433,25 +49439,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2114.3572
$ws.Range("I132").Value = 1255
$ws.Range("J132").Value = 4536.1816
$ws.Range("K132").Value = 3765
$ws.Range("L132").Value = 13608.5448
$ws.Range("M132").Value = -1235
$ws.Range("N132").Value = -18668.5448

# Hunk 28: sheet WVR, @@ -50574,25 +50580,25 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 56737.332
$ws.Range("J13").Value = 85006
$ws.Range("L13").Value = 85006
$ws.Range("N13").Value = -85286

# Hunk 29: sheet WVR, @@ -55168,22 +55174,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 561.2646999999999
$ws.Range("I107").Value = 251.26923
$ws.Range("K107").Value = 753.80769
$ws.Range("M107").Value = 1166.19231

# Hunk 30: sheet WVR, @@ -56393,25 +56399,25 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15898.333
$ws.Range("I132").Value = 2807.484
$ws.Range("J132").Value = 66625.375
$ws.Range("K132").Value = 8422.451999999999
$ws.Range("L132").Value = 199876.125
$ws.Range("M132").Value = -5892.451999999999
$ws.Range("N132").Value = -204936.125
